$d = $word.ActiveDocument

$replacements = @(
    @{old="63×53=3339"; new="35×35=1225"},
    @{old="42×33=1386"; new="51×82=4182"},
    @{old="81×24=1944"; new="18×82=1476"},
    @{old="89×11=979"; new="37×41=1517"},
    @{old="95×25=2375"; new="84×71=5964"},
    @{old="55×22=1210"; new="91×73=6643"},
    @{old="59×89=5251"; new="22×65=1430"},
    @{old="36×13=468"; new="45×70=3150"},
    @{old="69×96=6624"; new="47×92=4324"},
    @{old="96×79=7584"; new="94×83=7802"},
    @{old="71×93=6603"; new="98×55=5390"},
    @{old="28×13=364"; new="30×23=690"},
    @{old="19×23=437"; new="46×59=2714"},
    @{old="89×87=7743"; new="20×89=1780"},
    @{old="40×70=2800"; new="70×11=770"},
    @{old="81×69=5589"; new="37×93=3441"},
    @{old="94×26=2444"; new="27×25=675"},
    @{old="44×58=2552"; new="97×45=4365"},
    @{old="62×41=2542"; new="14×42=588"},
    @{old="16×29=464"; new="97×26=2522"},
    @{old="40×47=1880"; new="28×29=812"},
    @{old="63×77=4851"; new="16×40=640"},
    @{old="18×29=522"; new="80×46=3680"},
    @{old="48×15=720"; new="34×58=1972"},
    @{old="81×21=1701"; new="83×44=3652"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
